$d = $word.ActiveDocument
$g5text = "G5) A ride request must always be satisfied within a considerable short amount of time"
$null = $d.Content.Find.Execute($g5text, $true, $false, $false, $false, $false,
                         $true, 1, $false, $g5text, 2)

$p5 = $d.Paragraphs(5)
$g5start = $p5.Range.Start
$g5end = $g5start + $g5text.Length

$tailRange = $d.Range($g5end - 2, $g5end)
$tailRange.InsertAfter(" [within x minutes]")

Write-Output "done"
